$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - update values only
$ws.Range("B3").Value = 0.02821728926589555
$ws.Range("C3").Value = 0.02852433095172064
$ws.Range("D3").Value = 0.0294186858572493

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, update values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.02975141396014973
$ws.Range("C4").Value = 0.02965907731189859
$ws.Range("D4").Value = 0.02991523361315786

# Row 5: AdaBoostRegressor -> MLPRegressor, update values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.02593275180317702
$ws.Range("C5").Value = 0.02501107573167348
$ws.Range("D5").Value = 0.02423436760172389
